$wb = $excel.ActiveWorkbook
try {
  $p = $wb.BuiltinDocumentProperties.Item("Hyperlink base")
  Write-Output ("HyperlinkBase=" + $p.Value)
} catch { Write-Output "no hyperlink base" }
